$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows appended below the existing data (rows 15-27),
# mirroring rows 2-14 (same 13 matches) in a different order.
# Numeric-looking columns (G:K) are prefixed with a leading
# apostrophe so they stay text-typed, matching the existing rows.
$newRows = @(
    @(" Abu Dhabi", " October 25 2020", "Royals won by 8 wickets (with 10 balls remaining)", "Mumbai Indians", "Rajasthan Royals", "Ishan Kishan ", "'37", "'36", "'4", "'1", "'102.77"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Mumbai Indians", "Royal Challengers Bangalore", "Ishan Kishan ", "'25", "'19", "'3", "'1", "'131.57"),
    @(" Sharjah", " October 23 2020", "Mumbai won by 10 wickets (with 46 balls remaining)", "Mumbai Indians", "Chennai Super Kings", "Ishan Kishan ", "'68", "'37", "'6", "'5", "'183.78"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Ishan Kishan ", "'28", "'15", "'2", "'2", "'186.66"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Delhi Capitals", "Ishan Kishan ", "'55", "'30", "'4", "'3", "'183.33"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Mumbai Indians", "Kings XI Punjab", "Ishan Kishan ", "'28", "'32", "'1", "'1", "'87.50"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Mumbai Indians", "Sunrisers Hyderabad", "Ishan Kishan ", "'31", "'23", "'1", "'2", "'134.78"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Mumbai Indians", "Sunrisers Hyderabad", "Ishan Kishan ", "'33", "'30", "'1", "'2", "'110.00"),
    @(" Dubai (DSC)", " October 31 2020", "Mumbai won by 9 wickets (with 34 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Ishan Kishan ", "'72", "'47", "'8", "'3", "'153.19"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Mumbai Indians", "Delhi Capitals", "Ishan Kishan ", "'33", "'19", "'3", "'1", "'173.68"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Mumbai Indians", "Royal Challengers Bangalore", "Ishan Kishan ", "'99", "'58", "'2", "'9", "'170.68"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Rajasthan Royals", "Ishan Kishan ", "'0", "'1", "'0", "'0", "'0.00"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Mumbai Indians", "Kings XI Punjab", "Ishan Kishan ", "'7", "'7", "'1", "'0", "'100.00"),
)

$startRow = 15
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowValues = $newRows[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $rowValues[$j]
    }
}